$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plantilla Estudiantes")

# Convert CI column (D) values from text "CIxxxxxx" to plain numbers
$ws.Range("D2").Value = 123456
$ws.Range("D3").Value = 654321
$ws.Range("D4").Value = 789012

# Update the active selection to D2 (mirrors the saved selection state in the file)
$ws.Range("D2").Select()
